# Update the "Generate Report for Handback" timestamps across the
# Overview, zh-cn and de-de sheets of the handback status workbook.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview: Latest HO Xliff Generate Date
$wsOverview.Range("G2").Value = "2016-08-30 01:08:09"

# zh-cn: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H2").Value = "2016-08-30 01:07:59"
$wsZhCn.Range("K2").Value = "2016-08-30 01:08:28"

# de-de: Correspond Handoff Datetime / Correspond Handback DateTime
$wsDeDe.Range("H2").Value = "2016-08-30 01:08:09"
$wsDeDe.Range("K2").Value = "2016-08-30 01:08:35"
